$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent results for the 380 kV case (Case_4_15), rows 2-25
# correspond to columns B:O (14 lines in the OPF/loading result table).
$data = @(
    @(16.67922716492379, 9.632965862483871, 13.60841902117526, 14.26841639269425, 0, 36.27739680435264, 16.27977929898239, 25.41976137925523, 8.719340851021675, 0, 12.14883454313697, 16.78666016794496, 0, 25.76089114848177),
    @(16.16251391883767, 9.382459071776754, 13.63434335967346, 14.31617856532036, 0, 36.43403082400665, 16.34338254337887, 25.56296289311177, 8.722282103060175, 0, 12.14868545827624, 16.66279818476387, 0, 25.87232410170649),
    @(15.83747240677303, 9.224522586750117, 13.65267248625726, 14.34715210946347, 0, 36.54172794603139, 16.38524914224095, 25.65604171636899, 8.72431288902915, 0, 12.14996206210056, 16.58780344402676, 0, 25.9465351711187),
    @(15.70324218416679, 9.159194724790803, 13.66074795401442, 14.36018930332244, 0, 36.5884998758572, 16.40301772752783, 25.69526882283675, 8.725197194826812, 0, 12.15082748215871, 16.55753066780196, 0, 25.97823040405314),
    @(15.68085188583087, 9.148290778528791, 13.66212548504104, 14.36237922668747, 0, 36.5964401425693, 16.40601092602587, 25.70186080637947, 8.725347466346946, 0, 12.15099206568351, 16.55252194733416, 0, 25.98358109452657),
    @(15.83566905853966, 9.223645374090285, 13.65277894066718, 14.34732625109405, 0, 36.54234706717782, 16.38548591050109, 25.6565654949609, 8.724324585060838, 0, 12.14997233430729, 16.58739397868963, 0, 25.94695674181665),
    @(16.50278361054561, 9.547487411875315, 13.61685714478484, 14.28454346555867, 0, 36.32900481413476, 16.30112567146091, 25.46806850605321, 8.720308450543742, 0, 12.14849990075747, 16.74374673918236, 0, 25.79810984507776),
    @(17.7420131166431, 10.14700273091788, 13.56555422846359, 14.1744531675895, 0, 36.00263978910026, 16.15802209268181, 25.13925956538425, 8.71420800861441, 0, 12.15641482940944, 17.05774004640709, 0, 25.55228065479209),
    @(18.6017268790866, 10.56237934576975, 13.53953358460117, 14.10144858250993, 0, 35.81968877065586, 16.06649346806607, 24.92250690396454, 8.710796503333343, 0, 12.16874071203914, 17.29157402781728, 0, 25.39990054749265),
    @(18.98026083593876, 10.74527007167375, 13.53022924365725, 14.06993446623241, 0, 35.74894775358678, 16.02781000132628, 24.82927711691676, 8.709474503279088, 0, 12.17574418694946, 17.39834028012105, 0, 25.33674122009476),
    @(19.12168507547964, 10.81360755341761, 13.52706980150689, 14.05824379642459, 0, 35.72396643777396, 16.01358638185908, 24.79474518477258, 8.709006752177222, 0, 12.17859533848557, 17.43880168145564, 0, 25.31371282639365),
    @(19.09131382248775, 10.79893143912935, 13.52773406451008, 14.06075079532981, 0, 35.72926607805137, 16.01663078753873, 24.80214791642434, 8.709106032384573, 0, 12.17797246572274, 17.43008660324348, 0, 25.31863283167329),
    @(18.99193494325502, 10.75091091406522, 13.52996202313301, 14.06896780167097, 0, 35.74685625165287, 16.02663129690739, 24.82642068036667, 8.709435363658232, 0, 12.17597477301155, 17.40166858064252, 0, 25.3348288329099),
    @(18.93080954142727, 10.72137592261682, 13.53137409383292, 14.07403258165799, 0, 35.75786636668352, 16.03281224983533, 24.84138899927495, 8.709641362220568, 0, 12.17477700330838, 17.38426505225497, 0, 25.34486517113313),
    @(18.57672621821951, 10.55030082329704, 13.54019258658154, 14.10354216014252, 0, 35.82456408123161, 16.06908097638315, 24.92870775943597, 8.710887508494002, 0, 12.16831097309834, 17.28460236489736, 0, 25.40415235254508),
    @(18.35620765435367, 10.44376325726345, 13.54625094449733, 14.12207912082511, 0, 35.86868693092529, 16.09208728684326, 24.98365064666815, 8.711710722609203, 0, 12.16470073093027, 17.22354541921776, 0, 25.44210276342808),
    @(18.22819548170234, 10.38191720511101, 13.54997395387581, 14.13290077909214, 0, 35.89523966187879, 16.10559787495026, 25.01575800320142, 8.712205865769224, 0, 12.16275575619985, 17.18846567849853, 0, 25.46451038667673),
    @(18.18465448061897, 10.36088106500389, 13.55127545304402, 14.13659225929468, 0, 35.90443130781518, 16.1102200607524, 25.02671586176311, 8.712377237597924, 0, 12.16211986399696, 17.17659568282218, 0, 25.4721966695821),
    @(18.37980477435246, 10.45516355105311, 13.54558134983051, 14.12008930868951, 0, 35.86386834283948, 16.08960945530876, 24.97774955117868, 8.711620850587229, 0, 12.16507144544302, 17.23004122399241, 0, 25.43800287428145),
    @(19.02117783194444, 10.76504100713134, 13.52929774402586, 14.06654767831853, 0, 35.7416404674553, 16.02368236902496, 24.81927022714106, 8.709337740810513, 0, 12.17655615468925, 17.41001499362535, 0, 25.33004753270319),
    @(19.42912674600127, 10.96218695762201, 13.52077640115383, 14.03297137841146, 0, 35.67229359082819, 15.98307239706998, 24.72019549883029, 8.708037047175653, 0, 12.18522165197498, 17.52780924896246, 0, 25.26467370933587),
    @(19.21245803285764, 10.85747281267939, 13.52513045707789, 14.05076237110305, 0, 35.70833763783637, 16.00451994397355, 24.77266182330343, 8.708713799820806, 0, 12.18049120046888, 17.46493293306327, 0, 25.2990898844796),
    @(18.36914034855879, 10.45001133163671, 13.54588332605785, 14.12098839024888, 0, 35.86604313332699, 16.09072879845467, 24.98041581675983, 8.711661413600364, 0, 12.16490343847411, 17.22710439797556, 0, 25.43985459834678),
    @(17.41510766942957, 9.988994766445034, 13.5773832118366, 14.20284743709689, 0, 36.08100210072966, 16.1943453864107, 25.22384711972403, 8.715669506290373, 0, 12.1531245025677, 16.97214517940058, 0, 25.61383729168764)
)

$rowCount = $data.Count
$colCount = $data[0].Count
$arr = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$ws.Range("B2:O25").Value = $arr
